# Applies the Jan 5 2023 symbol-list refresh: updated prices/1h volumes for
# most rows, plus a corrected BKEXToken/KickToken row pair (rows 41-42 had
# their name/link/price/volume swapped & corrected).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns B/C hold plain text (coin name / link) so a direct .Value assignment
# is enough. Columns D/E hold numeric-looking text (price / percent change) that
# Excel would otherwise auto-convert to a number, so each of those cells is first
# switched to the Text number format ("@") before the new value is written - this
# is done per-cell, immediately before the write, so no other cell's formatting is
# touched.
$newValues = [ordered]@{
    'D2' = '257.03'
    'E2' = '-0.23%'
    'D3' = '27.04'
    'E3' = '-0.57%'
    'D4' = '4.658'
    'E4' = '-10.70%'
    'D5' = '0.05879'
    'E5' = '-0.72%'
    'E7' = '-0.46%'
    'D8' = '0.9502'
    'E8' = '-6.01%'
    'D9' = '0.1407'
    'E9' = '-0.56%'
    'D10' = '0.04096'
    'E10' = '14.91%'
    'D11' = '0.07088'
    'E11' = '-1.24%'
    'D12' = '0.03179'
    'E12' = '1.33%'
    'D13' = '0.09157'
    'E13' = '-0.88%'
    'D14' = '0.001539'
    'E14' = '-0.34%'
    'D15' = '0.0006051'
    'E15' = '-0.21%'
    'D16' = '0.006224'
    'E16' = '9.49%'
    'E17' = '1.02%'
    'D18' = '3.205'
    'E18' = '-1.90%'
    'E19' = '0.67%'
    'D20' = '0.3054'
    'E20' = '-2.91%'
    'E21' = '-0.47%'
    'D22' = '3.829'
    'E22' = '8.88%'
    'D23' = '0.04226'
    'E23' = '1.04%'
    'E24' = '0.15%'
    'D25' = '0.004296'
    'E25' = '-4.85%'
    'E26' = '0.11%'
    'D27' = '0.0001937'
    'E27' = '30.61%'
    'E40' = '0.50%'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D41' = '0.006211'
    'E41' = '10.04%'
    'B42' = 'BKEXToken'
    'C42' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D42' = '0.1101'
    'E42' = '-0.31%'
    'D43' = '0.002200'
    'E43' = '0.11%'
    'D44' = '0.01142'
    'E44' = '6.51%'
    'D45' = '0.00005471'
    'E45' = '0.76%'
    'E46' = '0.11%'
    'D47' = '0.07001'
    'E47' = '-35.76%'
    'D48' = '0.2328'
    'E48' = '10,322.98%'
    'D49' = '0.00002100'
    'E49' = '0.11%'
    'D50' = '0.0002000'
    'E50' = '0.11%'
}

foreach ($ref in $newValues.Keys) {
    $col = $ref.Substring(0, 1)
    if ($col -eq "D" -or $col -eq "E") {
        $ws.Range($ref).NumberFormat = "@"
    }
    $ws.Range($ref).Value = $newValues[$ref]
}
